$wb = $excel.ActiveWorkbook

# zh-cn sheet: update handoff/handback datetime for row 3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-18 12:24:43"
$wsZhCn.Range("G3").Value = "2016-01-18 12:25:26"

# de-de sheet: update handoff/handback datetime for row 3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-18 12:24:52"
$wsDeDe.Range("G3").Value = "2016-01-18 12:25:43"
